# Weekly price update: a new week's price observation for "Ají" (Inferno,
# Primera) is inserted as the new first data row (row 15), pushing all the
# existing observations down by one row (old row 15 -> 16, ..., old row 46 ->
# 47). The sheet's used dimension grows from A1:R46 to A1:R47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15, shifting row 15 and everything below it down.
$ws.Rows.Item(15).Insert(-4121)   # xlShiftDown

# Populate the newly inserted row 15 with this week's observation.
$ws.Cells.Item(15, 1).Value()  = 1
$ws.Cells.Item(15, 2).Value()  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value()  = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value()  = "2021-10-18"
$ws.Cells.Item(15, 5).Value()  = 15
$ws.Cells.Item(15, 6).Value()  = 100112021
$ws.Cells.Item(15, 7).Value()  = "Ají"
$ws.Cells.Item(15, 8).Value()  = "Inferno"
$ws.Cells.Item(15, 9).Value()  = "Primera"
$ws.Cells.Item(15, 10).Value() = 120
$ws.Cells.Item(15, 11).Value() = 36000
$ws.Cells.Item(15, 12).Value() = 37000
$ws.Cells.Item(15, 13).Value() = 36500
$ws.Cells.Item(15, 14).Value() = "$/caja 15 kilos"
$ws.Cells.Item(15, 15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item(15, 16).Value() = 2433
$ws.Cells.Item(15, 17).Value() = 15
$ws.Cells.Item(15, 18).Value() = "Hortaliza"
